$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-05 01:42:15"

for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
